$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2858
$ws.Range("J113").Value = 3376.5
$ws.Range("L113").Value = 3376.5
$ws.Range("N113").Value = -9884.5

$ws.Range("H131").Value = 3000
$ws.Range("I131").Value = 3000
$ws.Range("K131").Value = 9000
$ws.Range("M131").Value = -3960

$ws.Range("H132").Value = 4721.4707
$ws.Range("I132").Value = 2697.4644
$ws.Range("K132").Value = 8092.3932
$ws.Range("M132").Value = -5562.3932

$ws.Range("H138").Value = 1949.5354
$ws.Range("I138").Value = 806.3125
$ws.Range("J138").Value = 2169.9158
$ws.Range("K138").Value = 2418.9375
$ws.Range("L138").Value = 6509.7474
$ws.Range("M138").Value = 2721.0625
$ws.Range("N138").Value = -16789.7474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1669.2142
$ws.Range("I2").Value = 850
$ws.Range("K2").Value = 850
$ws.Range("M2").Value = -737

$ws.Range("H61").Value = 1633.2727
$ws.Range("I61").Value = 1496.2222
$ws.Range("K61").Value = 1496.2222
$ws.Range("M61").Value = -1284.2222

$ws.Range("H63").Value = 250001250
$ws.Range("I63").Value = 1660
$ws.Range("J63").Value = 1000000000
$ws.Range("K63").Value = 1660
$ws.Range("L63").Value = 1000000000
$ws.Range("M63").Value = -974
$ws.Range("N63").Value = -1000001372

$ws.Range("H66").Value = 250001250
$ws.Range("I66").Value = 1660
$ws.Range("J66").Value = 1000000000
$ws.Range("K66").Value = 8300
$ws.Range("L66").Value = 5000000000
$ws.Range("M66").Value = -4868
$ws.Range("N66").Value = -5000006864

$ws.Range("H74").Value = 619.7083
$ws.Range("I74").Value = 582.2895
$ws.Range("K74").Value = 582.2895
$ws.Range("M74").Value = 291.7105

$ws.Range("H77").Value = 619.7083
$ws.Range("I77").Value = 582.2895
$ws.Range("K77").Value = 2911.4475
$ws.Range("M77").Value = 1456.5525

$ws.Range("H116").Value = 1669.2142
$ws.Range("I116").Value = 850
$ws.Range("K116").Value = 850
$ws.Range("M116").Value = 1444

$ws.Range("H132").Value = 3347.2
$ws.Range("I132").Value = 3222.9524
$ws.Range("K132").Value = 9668.8572
$ws.Range("M132").Value = -7138.8572

$ws.Range("H136").Value = 1633.2727
$ws.Range("I136").Value = 1496.2222
$ws.Range("K136").Value = 4488.6666
$ws.Range("M136").Value = -1938.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1669.2142
$ws.Range("I3").Value = 850
$ws.Range("K3").Value = 850
$ws.Range("M3").Value = -736

$ws.Range("H134").Value = 13222.667
$ws.Range("I134").Value = 10350.75
$ws.Range("K134").Value = 31052.25
$ws.Range("M134").Value = -28517.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858260
$ws.Range("I16").Value = 142858260
$ws.Range("K16").Value = 142858260
$ws.Range("M16").Value = -142857973

$ws.Range("H31").Value = 664.914
$ws.Range("I31").Value = 587.8253999999999
$ws.Range("J31").Value = 826.8
$ws.Range("K31").Value = 587.8253999999999
$ws.Range("L31").Value = 826.8
$ws.Range("M31").Value = -292.8253999999999
$ws.Range("N31").Value = -1416.8

$ws.Range("H34").Value = 664.914
$ws.Range("I34").Value = 587.8253999999999
$ws.Range("J34").Value = 826.8
$ws.Range("K34").Value = 587.8253999999999
$ws.Range("L34").Value = 826.8
$ws.Range("M34").Value = -385.8253999999999
$ws.Range("N34").Value = -1230.8

$ws.Range("H94").Value = 1102.3334
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 1204.6666
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 1204.6666
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -2106.6666

$ws.Range("H113").Value = 142858260
$ws.Range("I113").Value = 142858260
$ws.Range("K113").Value = 142858260
$ws.Range("M113").Value = -142856090

$ws.Range("H132").Value = 11455.363
$ws.Range("I132").Value = 15430.143
$ws.Range("K132").Value = 46290.429
$ws.Range("M132").Value = -43760.429

$ws.Range("H133").Value = 63260.668
$ws.Range("J133").Value = 63260.668
$ws.Range("L133").Value = 63260.668
$ws.Range("N133").Value = -68320.66800000001

$ws.Range("H134").Value = 8773083
$ws.Range("I134").Value = 9524822
$ws.Range("K134").Value = 28574466
$ws.Range("M134").Value = -28571931

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 470
$ws.Range("I18").Value = 510
$ws.Range("K18").Value = 1530
$ws.Range("M18").Value = -1361

$ws.Range("H54").Value = 6005
$ws.Range("J54").Value = 6005
$ws.Range("L54").Value = 18015
$ws.Range("N54").Value = -19133

$ws.Range("H114").Value = 657.9545000000001
$ws.Range("J114").Value = 1032.4
$ws.Range("L114").Value = 3097.2
$ws.Range("N114").Value = -9605.200000000001

$ws.Range("H131").Value = 41668530
$ws.Range("I131").Value = 142857860
$ws.Range("J131").Value = 2335.4119
$ws.Range("K131").Value = 428573580
$ws.Range("L131").Value = 7006.2357
$ws.Range("M131").Value = -428568540
$ws.Range("N131").Value = -17086.2357

$ws.Range("H132").Value = 1964.2142
$ws.Range("I132").Value = 1375
$ws.Range("J132").Value = 2199.9
$ws.Range("K132").Value = 12375
$ws.Range("L132").Value = 19799.1
$ws.Range("M132").Value = -9845
$ws.Range("N132").Value = -24859.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6332.3335
$ws.Range("I43").Value = 4999
$ws.Range("K43").Value = 4999
$ws.Range("M43").Value = -4848

$ws.Range("H122").Value = 2365.9375
$ws.Range("I122").Value = 1319.625
$ws.Range("J122").Value = 3412.25
$ws.Range("K122").Value = 3958.875
$ws.Range("L122").Value = 10236.75
$ws.Range("M122").Value = -1508.875
$ws.Range("N122").Value = -15136.75

$ws.Range("H132").Value = 2511.375
$ws.Range("I132").Value = 2080.389
$ws.Range("J132").Value = 3804.3333
$ws.Range("K132").Value = 6241.167
$ws.Range("L132").Value = 11412.9999
$ws.Range("M132").Value = -3711.167
$ws.Range("N132").Value = -16472.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2254
$ws.Range("I7").Value = 2240.5
$ws.Range("J7").Value = 2267.5
$ws.Range("K7").Value = 2240.5
$ws.Range("L7").Value = 2267.5
$ws.Range("M7").Value = -2128.5
$ws.Range("N7").Value = -2491.5

$ws.Range("H22").Value = 1637.625
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1728.7142
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1728.7142
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2318.7142

$ws.Range("H27").Value = 1637.625
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1728.7142
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1728.7142
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1942.7142

$ws.Range("H40").Value = 2963
$ws.Range("I40").Value = 2890
$ws.Range("J40").Value = 2999.5
$ws.Range("K40").Value = 2890
$ws.Range("L40").Value = 2999.5
$ws.Range("M40").Value = -2754
$ws.Range("N40").Value = -3271.5

$ws.Range("H68").Value = 1530.9474
$ws.Range("I68").Value = 1266
$ws.Range("K68").Value = 1266
$ws.Range("M68").Value = -517

$ws.Range("H71").Value = 1530.9474
$ws.Range("I71").Value = 1266
$ws.Range("K71").Value = 6330
$ws.Range("M71").Value = -2586

$ws.Range("H126").Value = 2254
$ws.Range("I126").Value = 2240.5
$ws.Range("J126").Value = 2267.5
$ws.Range("K126").Value = 6721.5
$ws.Range("L126").Value = 6802.5
$ws.Range("M126").Value = -4251.5
$ws.Range("N126").Value = -11742.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 50006380
$ws.Range("I62").Value = 62505164
$ws.Range("J62").Value = 11250
$ws.Range("K62").Value = 62505164
$ws.Range("L62").Value = 11250
$ws.Range("M62").Value = -62504540
$ws.Range("N62").Value = -12498

$ws.Range("H65").Value = 50006380
$ws.Range("I65").Value = 62505164
$ws.Range("J65").Value = 11250
$ws.Range("K65").Value = 312525820
$ws.Range("L65").Value = 56250
$ws.Range("M65").Value = -312522700
$ws.Range("N65").Value = -62490

$ws.Range("H81").Value = 1420
$ws.Range("J81").Value = 4000
$ws.Range("L81").Value = 8000
$ws.Range("N81").Value = -10122

$ws.Range("H84").Value = 1420
$ws.Range("J84").Value = 4000
$ws.Range("L84").Value = 40000
$ws.Range("N84").Value = -50608

$ws.Range("H113").Value = 739.8333
$ws.Range("I113").Value = 319.77777
$ws.Range("K113").Value = 959.33331
$ws.Range("M113").Value = 1210.66669

$ws.Range("H126").Value = 50001616
$ws.Range("I126").Value = 125000960
$ws.Range("J126").Value = 2050.4167
$ws.Range("K126").Value = 375002880
$ws.Range("L126").Value = 6151.250100000001
$ws.Range("M126").Value = -375000410
$ws.Range("N126").Value = -11091.2501

$ws.Range("H132").Value = 2183.7407
$ws.Range("I132").Value = 1798.238
$ws.Range("J132").Value = 3533
$ws.Range("K132").Value = 5394.714
$ws.Range("L132").Value = 10599
$ws.Range("M132").Value = -2864.714
$ws.Range("N132").Value = -15659

$ws.Range("H141").Value = 37715
$ws.Range("J141").Value = 37715
$ws.Range("L141").Value = 37715
$ws.Range("N141").Value = -48075
